$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 / Spring 2022 / Summer 2022 block (rows 4-10) ---
# Spring 2022 column (C/D) - course list updated
$ws.Range("C4").Value = "KINS 2135"
$ws.Range("C5").Value = "KINS 2271"
$ws.Range("C6").Value = "KINS 2272"
$ws.Range("C7").Value = "KINS 2345"
$ws.Range("D7").Value = 1
$ws.Range("C8").Value = "KINS 2379"
$ws.Range("D8").Value = 1
$ws.Range("C9").Value = "KINS 3105"
$ws.Range("D9").Value = 3
$ws.Range("C10").Value = "KINS 3316"
$ws.Range("D10").Value = 1

# Summer 2022 column (E/F) - course list updated
$ws.Range("E4").Value = "CPSC 3121"
$ws.Range("E5").Value = "KINS 3126"
$ws.Range("F5").Value = 2
$ws.Range("E6").Value = "KINS 3165"
$ws.Range("F6").Value = 2
$ws.Range("E7").Value = "KINS 3256"
$ws.Range("F7").Value = 2

# Fall 2022 column (A/B) - new Psych/Kinesiology courses added, credits adjusted
$ws.Range("A5").Value = "PSYC 1101"
$ws.Range("B5").Value = 3
$ws.Range("A6").Value = "KINS 1105"
$ws.Range("A7").Value = "PSYC 1105"
$ws.Range("B7").Value = 2
$ws.Range("A8").Value = "KINS 1106"
$ws.Range("B8").Value = 2
$ws.Range("A9").Value = "KINS 2105"
$ws.Range("B9").Value = 2
$ws.Range("A10").Value = "CPSC 4000"
$ws.Range("B10").Value = 0

# --- Fall 2023 / Spring 2023 / Summer 2023 block (rows 13-19) ---
$ws.Range("A13").Value = "KINS 3107"
$ws.Range("C13").Value = "KINS 3235"
$ws.Range("E13").Value = "CPSC 4148"
$ws.Range("F13").Value = 3

$ws.Range("A14").Value = "DSCI 3111"
$ws.Range("C14").Value = "KINS 3255"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "CPSC 4698"
$ws.Range("F14").Value = 3

$ws.Range("A15").Value = "KINS 3127"
$ws.Range("C15").Value = "KINS 3257"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "CPSC 4899"
$ws.Range("F15").Value = 3

$ws.Range("A16").Value = "CPSC 3165"
$ws.Range("C16").Value = "KINS 3258"
$ws.Range("D16").Value = 2

$ws.Range("A17").Value = "KINS 3218"
$ws.Range("C17").Value = "KINS 3365"
$ws.Range("D17").Value = 1

# two extra Spring 2023 courses that didn't fit in the old layout
$ws.Range("C18").Value = "CPSC 3415"
$ws.Range("D18").Value = 1
$ws.Range("C19").Value = "CYBR 4125"
$ws.Range("D19").Value = 3

# --- Fall 2024 / Spring 2024 / Summer 2024 block (rows 22-26) ---
$ws.Range("A22").Value = "CPSC 4135"
$ws.Range("C22").Value = "CPSC 4176"
$ws.Range("D22").Value = 3

$ws.Range("A23").Value = "CYBR 4145"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = "CPSC 4205"
$ws.Range("D23").Value = 3

$ws.Range("A24").Value = "CPSC 4155"
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = "CYBR 4416"
$ws.Range("D24").Value = 1

$ws.Range("A25").Value = "CPSC 4157"
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = "CPSC 4555"
$ws.Range("D25").Value = 3

$ws.Range("A26").Value = "CPSC 4175"
$ws.Range("B26").Value = 3

# --- new Fall 2025 / Spring 2025 / Summer 2025 block (rows 30-38) ---
$ws.Range("A30").Value = "Fall 2025"
$ws.Range("B30").Value = "Credits"
$ws.Range("C30").Value = "Spring 2025"
$ws.Range("D30").Value = "Credits"
$ws.Range("E30").Value = "Summer 2025"
$ws.Range("F30").Value = "Credits"

$ws.Range("A38").Value = "Total"
$ws.Range("B38").Formula = "=SUM(B31:B37)"
$ws.Range("C38").Value = "Total"
$ws.Range("D38").Formula = "=SUM(D31:D37)"
$ws.Range("E38").Value = "Total"
$ws.Range("F38").Formula = "=SUM(F31:F37)"
